# Applies the "Add/update resource data for UWTranslationQuestions" edit.
#
# Summary of changes:
#  1. Remove the "License Information" (Heading2) paragraph.
#  2. In the paragraph that previously read:
#       "<bold>翻譯問題 (unfoldingWord)</bold> (Chinese (Traditional)) is based
#        on: unfoldingWord® Translation Questions, <link>unfoldingWord</link>,
#        2022, which is licensed under a <link>CC BY-SA 4.0 license</link>."
#     rename the bold run to "unfoldingWord® Translation Questions" and
#     replace everything that follows it (including the hyperlinks) with new
#     plain text describing the adapted-languages / license information.
#  3. Remove the paragraph "This PDF version is provided under the same
#     license." entirely (its content is now folded into the paragraph
#     above).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert `text` (never bold, inherits the zh_TW/zh_TW language that is
# already used throughout this document) at the collapsed position `pos`,
# producing its own distinct run without disturbing the adjoining runs.
# Returns the ending character position of the newly inserted text.
# ---------------------------------------------------------------------------
function Insert-PlainRun {
    param($pos, $text)

    $marker = "@@INSMARK@@"

    # Step 1: plant a small marker via straight Text assignment - this makes
    # the marker inherit the formatting (incl. w:lang) of whatever sits right
    # before the insertion point.
    $mk = $d.Range($pos, $pos)
    $mk.Text = $marker

    # Step 2: swap the marker for the real text through Find/Replace, forcing
    # the replacement to be non-bold; this keeps it from being re-absorbed
    # into a preceding bold run and gives it its own run.
    $fr = $d.Range(0, $d.Content.End)
    $fr.Find.ClearFormatting()
    $fr.Find.Replacement.ClearFormatting()
    $fr.Find.Replacement.Font.Bold = 0
    $null = $fr.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)

    # Step 3: locate the freshly-inserted text and neutralise the explicit
    # "Bold=False" toggle that step 2 introduced, without merging the run
    # into its neighbours.
    $loc = $d.Range($pos, $pos + $text.Length)
    $loc.Bold = 1
    $loc.Bold = 0

    return $pos + $text.Length
}

# ---------------------------------------------------------------------------
# Step 1: remove the "License Information" heading paragraph.
# ---------------------------------------------------------------------------
$licInfo = $d.Range(0, $d.Content.End)
$licInfo.Find.ClearFormatting()
$found = $licInfo.Find.Execute("License Information")
if ($found) {
    $para = $licInfo.Paragraphs(1)
    $para.Range.Delete()
}

# ---------------------------------------------------------------------------
# Step 2: locate the bold run "翻譯問題 (unfoldingWord)" and rename it.
# ---------------------------------------------------------------------------
$boldRng = $d.Range(0, $d.Content.End)
$boldRng.Find.ClearFormatting()
$boldRng.Find.Font.Bold = 1
$null = $boldRng.Find.Execute("翻譯問題 (unfoldingWord)")
$boldStart = $boldRng.Start

$boldRng.Text = "unfoldingWord® Translation Questions"
$boldEnd = $boldRng.End

# ---------------------------------------------------------------------------
# Step 3: wipe out everything from the end of the bold run through to (but
# excluding) the paragraph mark - this removes the old "(Chinese
# (Traditional)) is based on: ..." text and both hyperlinks.
# ---------------------------------------------------------------------------
$ownerPara = $d.Range($boldEnd, $boldEnd).Paragraphs(1)
$paraEnd = $ownerPara.Range.End
$tailRng = $d.Range($boldEnd, $paraEnd - 1)
$tailRng.Text = ""

# ---------------------------------------------------------------------------
# Step 4: remove the whole following paragraph ("This PDF version is
# provided under the same license.").
# ---------------------------------------------------------------------------
$nextParaRng = $d.Range($boldEnd, $boldEnd)
$nextParaRng.MoveStart(1, 1) | Out-Null
$nextPara = $nextParaRng.Paragraphs(1)
$nextPara.Range.Delete()

# ---------------------------------------------------------------------------
# Step 5: insert the new runs right after the bold run.
# ---------------------------------------------------------------------------
$pos = $boldEnd
$pos = Insert-PlainRun $pos " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "
$pos = Insert-PlainRun $pos "unfoldingWord® Translation Questions"
$pos = Insert-PlainRun $pos " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "
$pos = Insert-PlainRun $pos "unfoldingWord® Translation Questions"
$pos = Insert-PlainRun $pos " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

Write-Output "Done"
